# Add a new "Pulp Bleaching" unit process row (row 61) to the
# "Unit Processes" sheet, mirroring the structure of the existing rows.
#
# Columns: A=ID, B=meta-category, C=display name, D=product,
#          E=productType, F=varFile, G=varSheet, H=calcFile, I=calcSheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from the last populated row (59, which has
# the same style pattern needed for row 61: s=1 on most cells, s=2 on the
# "product" column D, and no explicit style on B) onto the new row.
$ws.Range("A59:I59").Copy()
$ws.Range("A61").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A61").Value = "paper_bleaching"
$ws.Range("B61").Value = "paper"
$ws.Range("C61").Value = "Pulp Bleaching"
$ws.Range("D61").Value = "bleached pulp"
$ws.Range("E61").Value = "outflow"
$ws.Range("F61").Value = "data/paper/paper_var.xlsx"
$ws.Range("G61").Value = "Bleach"
$ws.Range("H61").Value = "data/paper/paper_calc.xlsx"
$ws.Range("I61").Value = "Bleach"
